$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-07-24 Monday" "2023-07-25 Tuesday"

Replace-Text "11÷8=" "27÷5="
Replace-Text "29÷4=" "87÷7="
Replace-Text "16÷3=" "21÷2="
Replace-Text "32÷4=" "76÷9="
Replace-Text "84÷6=" "11÷3="

Replace-Text "58÷8=" "92÷7="
Replace-Text "22÷8=" "76÷9="
Replace-Text "74÷4=" "37÷3="
Replace-Text "15÷6=" "18÷5="
Replace-Text "28÷4=" "70÷9="

Replace-Text "33÷6=" "48÷6="
Replace-Text "13÷7=" "59÷2="
Replace-Text "93÷2=" "98÷9="
Replace-Text "98÷3=" "86÷3="
Replace-Text "96÷5=" "48÷9="

Replace-Text "84÷8=" "37÷6="
Replace-Text "72÷6=" "21÷8="
Replace-Text "56÷5=" "51÷2="
Replace-Text "61÷7=" "16÷3="
Replace-Text "33÷2=" "36÷2="

Replace-Text "79÷9=" "90÷8="
Replace-Text "14÷7=" "56÷3="
Replace-Text "84÷5=" "99÷7="
Replace-Text "74÷8=" "90÷4="
Replace-Text "16÷2=" "95÷4="

Write-Output "Done"
